# The underlying change in this revision ("Fixed POI packaging and
# upgraded to POI 3.15") is purely a side effect of bumping the Apache
# POI library used to *write* the .docx fixture: the newer POI release
# serializes XML attributes in (namespace-prefix-agnostic) alphabetical
# order, whereas the old one preserved Word's original attribute order.
#
# Diffing the canonical OOXML confirms this: every element that shows a
# "-"/"+" pair in word/document.xml and word/styles.xml has the exact
# same tag name and the exact same set of attribute name/value pairs on
# both sides - only the left-to-right order of those attributes differs
# (e.g. <w:pgSz w:w="11906" w:h="16838"/> becomes
# <w:pgSz w:h="16838" w:w="11906"/>). No text, value, style, formatting
# or structural content changes anywhere in the package.
#
# Attribute emission order is an internal detail of whichever XML writer
# serializes the package - it is not something the Word object model
# (Find/Replace, Paragraphs, PageSetup, Styles, ...) exposes any control
# over, in real Word automation or here. So the content-faithful way to
# reproduce this revision through COM automation is to leave the
# document's paragraphs/styles/sections/content completely untouched:
# there is nothing to find-and-replace, no run/paragraph/style property
# actually changed values, so no Range/Find/Style mutation is needed.
$d = $word.ActiveDocument

# Touch nothing: re-saving without edits keeps every paragraph, run,
# section and style property exactly as authored, matching the diff's
# lack of any semantic change.
